# Data for the new "Work-Items2023" sheet.
# Columns: A=WIID, B=Description, C=Type, D=Status, E=Date
$workItems = @(
    @('WIID','Description','Type','Status','Date'),
    @('92928123','Verify Account Position','WI1','Open','2023-05-25'),
    @('92928061','Process Vendor Invoice','WI3','Open','2023-01-19'),
    @('92928134','Verify Account Position','WI1','Open','2023-05-10'),
    @('92928077','Research Client Check Copy','WI2','Open','2023-04-03'),
    @('92928072','Research Client Check Copy','WI2','Open','2023-05-26'),
    @('92928034','Calculate Client Security Hash','WI5','Open','2023-01-08'),
    @('92928127','Verify Account Position','WI1','Open','2023-01-18'),
    @('92928145','Verify Account Position','WI1','Open','2023-04-27'),
    @('92928068','Research Client Check Copy','WI2','Open','2023-04-03'),
    @('92928024','Calculate Client Security Hash','WI5','Open','2023-03-26'),
    @('92928084','Research Client Check Copy','WI2','Open','2023-02-12'),
    @('92928064','Process Vendor Invoice','WI3','Open','2023-04-28')
)

$wb = $excel.ActiveWorkbook

# --- Rename the original sheet to "Master" ---
$master = $wb.ActiveSheet
$master.Name = "Master"

# --- Add the new sheet right after "Master" and name it ---
$wi = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $master)
$wi.Name = "Work-Items2023"

# --- Populate "Work-Items2023" with the work item table ---
for ($r = 0; $r -lt $workItems.Length; $r++) {
    $row = $workItems[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $cell = $wi.Cells.Item($r + 1, $c + 1)
        $val = $row[$c]
        # Column A (WIID) and column E (Date) are all-digit / date-looking
        # text values that must stay as text, not be auto-converted to a
        # number or a date serial.
        if ($c -eq 0 -or $c -eq 4) {
            $cell.Value = "'" + $val
        } else {
            $cell.Value = $val
        }
    }
}

# Column B ("Description") is the widest column on this sheet.
$wi.Columns.Item(2).ColumnWidth = 27.71

# --- View / selection state to match the edited workbook ---
# "Master" keeps cell A2 selected and is no longer the active tab.
$master.Range("A2").Select()

# "Work-Items2023" becomes the active sheet with G17 selected.
$wi.Activate()
$wi.Range("G17").Select()
